$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update F2, F5, F6
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9987
$ws1.Range("F5").Value = 599
$ws1.Range("F6").Value = 476

# Sheet "全部类型" (sheet4): update F2, F5, F7
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9987
$ws4.Range("F5").Value = 599
$ws4.Range("F7").Value = 476
